$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 94949
$ws.Range("B2").Value = "Nicolas da Cruz"
$ws.Range("C2").Value = "P&D"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 45105
$ws.Range("G2").Value = 2843.11

# Row 3
$ws.Range("A3").Value = 44519
$ws.Range("B3").Value = "Theo Lima"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45096
$ws.Range("G3").Value = 4509.77

# Row 4
$ws.Range("A4").Value = 92271
$ws.Range("B4").Value = "Diogo Carvalho"
$ws.Range("C4").Value = "Juridico"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 8016.64

# Row 5
$ws.Range("A5").Value = 82343
$ws.Range("B5").Value = "Luiza da Cunha"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Viagem de negocios"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45096
$ws.Range("G5").Value = 9051.77

# Row 6
$ws.Range("A6").Value = 20209
$ws.Range("B6").Value = "Gabriel da Cruz"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Viagem de negocios"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45078
$ws.Range("G6").Value = 6698.88

# Row 7
$ws.Range("A7").Value = 92477
$ws.Range("B7").Value = "Maria Júlia Castro"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45079
$ws.Range("G7").Value = 5014.56

# Row 8
$ws.Range("A8").Value = 93705
$ws.Range("B8").Value = "Emilly Pastor"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Consulta medica"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45087
$ws.Range("G8").Value = 9181.95

# Row 9
$ws.Range("A9").Value = 99241
$ws.Range("B9").Value = "Rebeca da Costa"
$ws.Range("C9").Value = "Juridico"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45099
$ws.Range("G9").Value = 8836.65

# Row 10
$ws.Range("A10").Value = 85438
$ws.Range("B10").Value = "Laís Vasconcelos"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Doenca"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45084
$ws.Range("G10").Value = 4257.46

# Row 11
$ws.Range("A11").Value = 22564
$ws.Range("B11").Value = "Luara Araújo"
$ws.Range("C11").Value = "Operacoes"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45092
$ws.Range("G11").Value = 3225.29

$wb.Save()
